# Add "0" as a new, first option for the barcode_offset pick-list.
$wb = $excel.ActiveWorkbook

$offsetSheet = $wb.Worksheets.Item("barcode_offset")

# Insert a new row at the top of the list and give it the value "0".
# The leading apostrophe forces Excel to store it as text (matching the
# other entries in this list, e.g. "8"), and ClearFormats() drops the
# "number stored as text" quote-prefix styling that would otherwise be
# picked up, so the cell ends up as a plain shared-string cell.
$offsetSheet.Rows.Item(1).Insert()
$newCell = $offsetSheet.Cells.Item(1, 1)
$newCell.Value = "'0"
$newCell.ClearFormats()

# The validation list on the main sheet points at a fixed range
# ('barcode_offset'!$A$1:$A$5); grow it by one row so the new option is
# selectable, while preserving the existing validation behaviour.
$mainSheet = $wb.Worksheets.Item("RNAseq")
$offsetValidation = $mainSheet.Range("O2:O1001").Validation
$offsetValidation.Formula1 = "='barcode_offset'!`$A`$1:`$A`$6"
$offsetValidation.IgnoreBlank = $true
$offsetValidation.ShowError = $true

# The template's own ".metadata" sheet stamps the export time; bump it to
# match the re-export that carried this edit.
$metaSheet = $wb.Worksheets.Item(".metadata")
$metaSheet.Cells.Item(2, 3).Value = "2023-10-31T14:33:40-07:00"
